$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U3").Value = 1.92
$ws.Range("V3").Value = 1.77
$ws.Range("BD4").Value = 126
$ws.Range("V4").Value = 1.63
$ws.Range("BD7").Value = 126
$ws.Range("Q7").Value = 1.95
$ws.Range("R7").Value = 1.9
$ws.Range("M11").Value = 1.02
$ws.Range("O11").Value = 1.13
$ws.Range("M12").Value = 1.03
$ws.Range("O12").Value = 1.17
$ws.Range("Q12").Value = 1.7
$ws.Range("R12").Value = 2.1
$ws.Range("M13").Value = 1.01
$ws.Range("O13").Value = 1.08
$ws.Range("AD15").Value = 8
$ws.Range("K15").Value = 2.4
$ws.Range("O15").Value = 1.2
$ws.Range("P15").Value = 4.33
$ws.Range("Q15").Value = 1.67
$ws.Range("R15").Value = 2.15
$ws.Range("U15").Value = 1.62
$ws.Range("V15").Value = 2.2
$ws.Range("AA16").Value = 13
$ws.Range("AD16").Value = 8
$ws.Range("AT16").Value = 3
$ws.Range("AX16").Value = 34
$ws.Range("K16").Value = 2.3
$ws.Range("N16").Value = 13
$ws.Range("S16").Value = 1.36
$ws.Range("T16").Value = 3
$ws.Range("U16").Value = 1.91
$ws.Range("V16").Value = 1.8
$ws.Range("X16").Value = 7
$ws.Range("AF17").Value = 26
$ws.Range("AH17").Value = 17
$ws.Range("AI17").Value = 19
$ws.Range("AL17").Value = 19
$ws.Range("AM17").Value = 21
$ws.Range("AP17").Value = 15
$ws.Range("AQ17").Value = 34
$ws.Range("AW17").Value = 5.5
$ws.Range("G17").Value = 2.3
$ws.Range("I17").Value = 2.7
$ws.Range("J17").Value = 2.75
$ws.Range("L17").Value = 3.1
$ws.Range("N17").Value = 23
$ws.Range("O17").Value = 1.11
$ws.Range("U17").Value = 1.37
$ws.Range("Z17").Value = 23
$ws.Range("M18").Value = 1.05
$ws.Range("O18").Value = 1.29
$ws.Range("U18").Value = 1.69
$ws.Range("AB19").Value = 23
$ws.Range("AD19").Value = 8
$ws.Range("AU19").Value = 8
$ws.Range("AY19").Value = 29
$ws.Range("AZ19").Value = 81
$ws.Range("G19").Value = 1.67
$ws.Range("H19").Value = 4
$ws.Range("I19").Value = 4.5
$ws.Range("J19").Value = 2.25
$ws.Range("M19").Value = 1.03
$ws.Range("O19").Value = 1.2
$ws.Range("U19").Value = 1.63
$ws.Range("AB20").Value = 34
$ws.Range("AE20").Value = 15
$ws.Range("AG20").Value = 126
$ws.Range("AH20").Value = 10
$ws.Range("AI20").Value = 9
$ws.Range("AK20").Value = 12
$ws.Range("AM20").Value = 19
$ws.Range("AN20").Value = 7.5
$ws.Range("AO20").Value = 26
$ws.Range("AP20").Value = 26
$ws.Range("AQ20").Value = 81
$ws.Range("AR20").Value = 81
$ws.Range("AU20").Value = 7.5
$ws.Range("AX20").Value = 7.5
$ws.Range("G20").Value = 5.25
$ws.Range("H20").Value = 4.5
$ws.Range("I20").Value = 1.48
$ws.Range("L20").Value = 2
$ws.Range("U20").Value = 1.54
$ws.Range("V20").Value = 2.25
$ws.Range("Z20").Value = 51
$ws.Range("AB21").Value = 17
$ws.Range("AD21").Value = 11
$ws.Range("AK21").Value = 67
$ws.Range("AM21").Value = 34
$ws.Range("AW21").Value = 8.5
$ws.Range("AZ21").Value = 81
$ws.Range("G21").Value = 1.38
$ws.Range("H21").Value = 5.25
$ws.Range("I21").Value = 6.5
$ws.Range("U21").Value = 1.47
$ws.Range("V21").Value = 2.5
$ws.Range("U22").Value = 1.54
$ws.Range("U24").Value = 1.91
$ws.Range("V24").Value = 1.8
$ws.Range("AA25").Value = 13
$ws.Range("AD25").Value = 7.5
$ws.Range("AE25").Value = 17
$ws.Range("AH25").Value = 13
$ws.Range("AI25").Value = 23
$ws.Range("AJ25").Value = 15
$ws.Range("AK25").Value = 51
$ws.Range("AO25").Value = 9
$ws.Range("AP25").Value = 19
$ws.Range("AQ25").Value = 29
$ws.Range("AS25").Value = 126
$ws.Range("AT25").Value = 3
$ws.Range("AW25").Value = 6.5
$ws.Range("AX25").Value = 23
$ws.Range("G25").Value = 1.7
$ws.Range("H25").Value = 3.8
$ws.Range("I25").Value = 4.5
$ws.Range("J25").Value = 2.3
$ws.Range("K25").Value = 2.25
$ws.Range("L25").Value = 5
$ws.Range("M25").Value = 1.04
$ws.Range("N25").Value = 13
$ws.Range("Q25").Value = 1.88
$ws.Range("R25").Value = 1.98
$ws.Range("S25").Value = 1.36
$ws.Range("T25").Value = 3
$ws.Range("U25").Value = 1.8
$ws.Range("V25").Value = 1.95
$ws.Range("W25").Value = 7
$ws.Range("X25").Value = 8
$ws.Range("Z25").Value = 13
$ws.Range("M28").Value = 1.05
$ws.Range("N28").Value = 11
$ws.Range("G30").Value = 2.1
$ws.Range("U30").Value = 1.8
$ws.Range("V30").Value = 1.91
$ws.Range("AB32").Value = 23
$ws.Range("AC32").Value = 13
$ws.Range("AD32").Value = 7
$ws.Range("AG32").Value = 151
$ws.Range("AH32").Value = 12
$ws.Range("AM32").Value = 29
$ws.Range("AP32").Value = 19
$ws.Range("AQ32").Value = 41
$ws.Range("AS32").Value = 126
$ws.Range("AT32").Value = 3.25
$ws.Range("AU32").Value = 7.5
$ws.Range("AY32").Value = 23
$ws.Range("BB32").Value = 151
$ws.Range("BC32").Value = 451
$ws.Range("G32").Value = 2.1
$ws.Range("H32").Value = 3.5
$ws.Range("J32").Value = 2.75
$ws.Range("K32").Value = 2.3
$ws.Range("L32").Value = 3.75
$ws.Range("M32").Value = 1.03
$ws.Range("N32").Value = 13
$ws.Range("O32").Value = 1.19
$ws.Range("P32").Value = 4
$ws.Range("Q32").Value = 1.73
$ws.Range("R32").Value = 2.08
$ws.Range("S32").Value = 1.33
$ws.Range("T32").Value = 3.25
$ws.Range("U32").Value = 1.62
$ws.Range("V32").Value = 2.2
$ws.Range("W32").Value = 9.5
$ws.Range("X32").Value = 11
$ws.Range("AI33").Value = 12
$ws.Range("AK33").Value = 26
$ws.Range("BD33").Value = 151
$ws.Range("G33").Value = 2.88
$ws.Range("M33").Value = 1.05
$ws.Range("N33").Value = 8
$ws.Range("O33").Value = 1.37
$ws.Range("U33").Value = 1.95
$ws.Range("V33").Value = 1.8
$ws.Range("AB34").Value = 26
$ws.Range("AC34").Value = 12
$ws.Range("AG34").Value = 301
$ws.Range("AO34").Value = 8.5
$ws.Range("AT34").Value = 3
$ws.Range("BD34").Value = 126
$ws.Range("G34").Value = 1.62
$ws.Range("M34").Value = 1.03
$ws.Range("O34").Value = 1.22
$ws.Range("Q34").Value = 1.88
$ws.Range("R34").Value = 1.98
$ws.Range("S34").Value = 1.36
$ws.Range("T34").Value = 3
$ws.Range("U34").Value = 1.91
$ws.Range("V34").Value = 1.91
$ws.Range("W34").Value = 7
$ws.Range("AH35").Value = 8.5
$ws.Range("AL35").Value = 26
$ws.Range("G35").Value = 2.4
$ws.Range("U35").Value = 1.83
$ws.Range("V35").Value = 1.83
$ws.Range("W35").Value = 7.5
$ws.Range("Y35").Value = 10
$ws.Range("Z35").Value = 23
$ws.Range("AG36").Value = 251
$ws.Range("AX36").Value = 19
$ws.Range("BB36").Value = 201
$ws.Range("G36").Value = 2.1
$ws.Range("H36").Value = 3.3
$ws.Range("I36").Value = 3.2
$ws.Range("U36").Value = 1.8
$ws.Range("V36").Value = 1.91
$ws.Range("X36").Value = 10
$ws.Range("Z36").Value = 19
$ws.Range("M38").Value = 1.05
$ws.Range("O38").Value = 1.41
$ws.Range("P38").Value = 2.62
$ws.Range("M39").Value = 1.03
$ws.Range("O39").Value = 1.25
